# Insert a new weekly price record as row 215 on the "Hortaliza, Femacal de
# La Calera - Haba" sheet. This shifts all the existing rows from 215..249
# down to 216..250 (data unchanged), and populates the newly inserted row
# 215 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 215, pushing rows 215:249 down to 216:250.
$ws.Rows.Item(215).EntireRow.Insert()

# Fill in the new row 215 with the new record.
$ws.Range("A215").Value = 3
$ws.Range("B215").Value = "Femacal de La Calera"
$ws.Range("C215").Value = "Coquimbo"
$ws.Range("D215").Value = 45077
$ws.Range("E215").Value = 5
$ws.Range("F215").Value = 100112026
$ws.Range("G215").Value = "Haba"
$ws.Range("H215").Value = "Sin especificar"
$ws.Range("I215").Value = "Primera"
$ws.Range("J215").Value = 78
$ws.Range("K215").Value = 20000
$ws.Range("L215").Value = 21000
$ws.Range("M215").Value = 20487
$ws.Range("N215").Value = "`$/saco 25 kilos"
$ws.Range("O215").Value = "Provincia de Limar$([char]0x00ED)"
$ws.Range("P215").Value = 819
$ws.Range("Q215").Value = 25
$ws.Range("R215").Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D.
$ws.Range("D215").NumberFormat = $ws.Range("D216").NumberFormat
